$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary updates -------------------------------------------------
$ws.Range("E11").Value = 269662
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 9

# Swap the "Novedad de Retiro" / "Novedad de Ingreso" column headers
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# --- Preserve the special "last data row" border/format -------------------
# Row 26 (old) -> row 32 in the original layout carries a distinct bottom
# border style. Copy that formatting onto what will become the new last
# data row (row 26) before the intervening rows are removed.
$ws.Range("B32:J32").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the rows that are no longer part of the account statement -----
# (old rows 27-30 = JOSE DANIEL RODRIGUEZ BARRETO periods, old rows 31-32 =
# duplicate JESSICA ANDREINA MARTINEZ MOROS periods now folded into rows
# 23-24). Deleting these 6 rows shifts the signature-block rows (old 37-38)
# up to rows 31-32, matching the new, shorter table.
$ws.Range("27:32").EntireRow.Delete()

# --- Rewrite the worker/period data table (rows 16-26) ---------------------
$data = @(
    @("CC", "1050962594", "KARL HELLS CHAMORRO RAMOS", "1911", 37800, 945000),
    @("CC", "20238215", "ANDERSON YAFETH BLANCO HERNANDEZ", "1911", 14354, 828116),
    @("CC", "1050962594", "KARL HELLS CHAMORRO RAMOS", "1912", 37800, 945000),
    @("CC", "1050962594", "KARL HELLS CHAMORRO RAMOS", "2001", 37800, 945000),
    @("CC", "1143341319", "KITYAN MARCELA CASTRO PALACIO", "2110", 21804, 1000000),
    @("CC", "1128056330", "KARELIS MARIA MOLINA TORRES", "2201", 18666, 1000000),
    @("PE", "6207678", "VERONICA VALENTINA BOLIVAR REMOLINA", "2207", 1333, 1000000),
    @("CC", "1143373218", "JESSICA ANDREINA MARTINEZ MOROS", "2207", 52000, 1300000),
    @("CC", "1143373218", "JESSICA ANDREINA MARTINEZ MOROS", "2208", 15600, 1300000),
    @("CC", "1128062875", "LUIS ENRIQUE DRITT TRUJILLO", "2210", 4665, 3498894),
    @("CC", "1001978144", "VALENTINA DE AVILA JULIAO", "2309", 27840, 1160000)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
